$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '29.813.47'
$ws.Range('E2').Value = '  -1.17%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.899.35'
$ws.Range('E3').Value = '  -0.77%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  +0.20%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '0.7653'
$ws.Range('E5').Value = '  +3.31%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '240.39'
$ws.Range('E6').Value = '  -1.45%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3043'
$ws.Range('E8').Value = '  -2.66%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '25.28'
$ws.Range('E9').Value = '  -5.09%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.06839'
$ws.Range('E10').Value = '  -1.77%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07982'
$ws.Range('E11').Value = '  +0.01%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '1.897.35'
$ws.Range('E12').Value = '  -0.42%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.7343'
$ws.Range('E13').Value = '  -5.99%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '5.163'
$ws.Range('E14').Value = '  -2.29%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '91.04'
$ws.Range('E15').Value = '  -1.43%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '29.827.23'
$ws.Range('E16').Value = '  -1.16%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '13.71'
$ws.Range('E17').Value = '  -4.46%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '5.879'
$ws.Range('E18').Value = '  -0.37%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '244.87'
$ws.Range('E19').Value = '  +1.05%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '0.000007694'
$ws.Range('E20').Value = '  -1.69%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '1.001'
$ws.Range('E21').Value = '  +0.14%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '2.135.52'
$ws.Range('E22').Value = '  -0.95%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '1.002'
$ws.Range('E23').Value = '  +0.27%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '6.872'
$ws.Range('E24').Value = '  -2.17%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '166.93'
$ws.Range('E25').Value = '  -0.68%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '9.230'
$ws.Range('E26').Value = '  -1.87%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '18.65'
$ws.Range('E27').Value = '  -2.45%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '0.1280'
$ws.Range('E28').Value = '  -0.32%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '2.024'
$ws.Range('E29').Value = '  -2.54%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.400'
$ws.Range('E30').Value = '  +3.71%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.511'
$ws.Range('E31').Value = '  -2.27%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '4.264'
$ws.Range('E32').Value = '  -1.70%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '4.062'
$ws.Range('E33').Value = '  -1.07%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.05252'
$ws.Range('E34').Value = '  +1.61%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.240'
$ws.Range('E35').Value = '  -5.01%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.7250'
$ws.Range('E36').Value = '  -2.95%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.719'
$ws.Range('E37').Value = '  -0.17%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.01908'
$ws.Range('E38').Value = '  -1.96%  '
$ws.Range('E39').Value = '  -0.73%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '6.200'
$ws.Range('E40').Value = '  -2.47%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.4399'
$ws.Range('E41').Value = '  -2.36%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '71.89'
$ws.Range('E42').Value = '  -4.22%  '
$ws.Range('E43').Value = '  +0.10%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.8354'
$ws.Range('E44').Value = '  -0.34%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '1.878'
$ws.Range('E45').Value = '  -4.02%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '7.577'
$ws.Range('E46').Value = '  -3.91%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '99.85'
$ws.Range('E47').Value = '  -1.66%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '9.724'
$ws.Range('E48').Value = '  -1.94%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.045.23'
$ws.Range('E49').Value = '  +0.38%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '36.13'
$ws.Range('E50').Value = '  -3.47%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.05919'
$ws.Range('E51').Value = '  -1.26%  '
